$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58 - constant columns (new row)
$ws.Range("A58").Value = 7
$ws.Range("B58").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C58").Value = 'Ñuble'
$ws.Range("E58").Value = 16
$ws.Range("F58").Value = 'Fruta'
$ws.Range("G58").Value = 100104
$ws.Range("H58").Value = 'Frutos de pepita'
$ws.Range("I58").Value = 100104003
$ws.Range("J58").Value = 'Membrillo'
$ws.Range("K58").Value = 'Champion'
$ws.Range("D58").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 59 - constant columns (new row)
$ws.Range("A59").Value = 7
$ws.Range("B59").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C59").Value = 'Ñuble'
$ws.Range("E59").Value = 16
$ws.Range("F59").Value = 'Fruta'
$ws.Range("G59").Value = 100104
$ws.Range("H59").Value = 'Frutos de pepita'
$ws.Range("I59").Value = 100104003
$ws.Range("J59").Value = 'Membrillo'
$ws.Range("K59").Value = 'Champion'
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 30
$ws.Range("D30").Value = 45119
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("Q30").Value = '$/bandeja 18 kilos granel'
$ws.Range("R30").Value = 'Región de O''Higgins'
$ws.Range("S30").Value = 556
$ws.Range("T30").Value = 18

# Row 31
$ws.Range("D31").Value = 45119
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 8000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 8000
$ws.Range("Q31").Value = '$/bandeja 18 kilos granel'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 444
$ws.Range("T31").Value = 18

# Row 32
$ws.Range("D32").Value = 45049
$ws.Range("L32").Value = 'Especial'
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 13000
$ws.Range("O32").Value = 13000
$ws.Range("P32").Value = 13000
$ws.Range("Q32").Value = '$/caja 18 kilos empedrada'
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 722
$ws.Range("T32").Value = 18

# Row 33
$ws.Range("D33").Value = 45049
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 12000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 12000
$ws.Range("Q33").Value = '$/caja 18 kilos empedrada'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("S33").Value = 667
$ws.Range("T33").Value = 18

# Row 34
$ws.Range("D34").Value = 45090
$ws.Range("L34").Value = 'Especial'
$ws.Range("M34").Value = 80
$ws.Range("N34").Value = 11000
$ws.Range("O34").Value = 11000
$ws.Range("P34").Value = 11000
$ws.Range("Q34").Value = '$/caja 18 kilos empedrada'
$ws.Range("R34").Value = 'Región del Maule'
$ws.Range("S34").Value = 611
$ws.Range("T34").Value = 18

# Row 35
$ws.Range("D35").Value = 45090
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 10000
$ws.Range("O35").Value = 10000
$ws.Range("P35").Value = 10000
$ws.Range("Q35").Value = '$/caja 18 kilos empedrada'
$ws.Range("R35").Value = 'Región del Maule'
$ws.Range("S35").Value = 556
$ws.Range("T35").Value = 18

# Row 36
$ws.Range("D36").Value = 45090
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 40
$ws.Range("N36").Value = 8000
$ws.Range("O36").Value = 8000
$ws.Range("P36").Value = 8000
$ws.Range("Q36").Value = '$/caja 18 kilos empedrada'
$ws.Range("R36").Value = 'Región del Maule'
$ws.Range("S36").Value = 444
$ws.Range("T36").Value = 18

# Row 37
$ws.Range("D37").Value = 44699
$ws.Range("L37").Value = 'Especial'
$ws.Range("M37").Value = 60
$ws.Range("N37").Value = 13000
$ws.Range("O37").Value = 13000
$ws.Range("P37").Value = 13000
$ws.Range("Q37").Value = '$/caja 15 kilos granel'
$ws.Range("R37").Value = 'Provincia de Curicó'
$ws.Range("S37").Value = 867
$ws.Range("T37").Value = 15

# Row 38
$ws.Range("D38").Value = 44699
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 120
$ws.Range("N38").Value = 11000
$ws.Range("O38").Value = 12000
$ws.Range("P38").Value = 11500
$ws.Range("Q38").Value = '$/caja 15 kilos granel'
$ws.Range("R38").Value = 'Provincia de Curicó'
$ws.Range("S38").Value = 767
$ws.Range("T38").Value = 15

# Row 39
$ws.Range("D39").Value = 45085
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 50
$ws.Range("N39").Value = 10000
$ws.Range("O39").Value = 10000
$ws.Range("P39").Value = 10000
$ws.Range("Q39").Value = '$/caja 18 kilos empedrada'
$ws.Range("R39").Value = 'Región del Maule'
$ws.Range("S39").Value = 556
$ws.Range("T39").Value = 18

# Row 40
$ws.Range("D40").Value = 45062
$ws.Range("L40").Value = 'Especial'
$ws.Range("M40").Value = 50
$ws.Range("N40").Value = 13000
$ws.Range("O40").Value = 13000
$ws.Range("P40").Value = 13000
$ws.Range("Q40").Value = '$/caja 18 kilos empedrada'
$ws.Range("R40").Value = 'Región de O''Higgins'
$ws.Range("S40").Value = 722
$ws.Range("T40").Value = 18

# Row 41
$ws.Range("D41").Value = 45062
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 50
$ws.Range("N41").Value = 12000
$ws.Range("O41").Value = 12000
$ws.Range("P41").Value = 12000
$ws.Range("Q41").Value = '$/caja 18 kilos empedrada'
$ws.Range("R41").Value = 'Región de O''Higgins'
$ws.Range("S41").Value = 667
$ws.Range("T41").Value = 18

# Row 42
$ws.Range("D42").Value = 45079
$ws.Range("L42").Value = 'Especial'
$ws.Range("M42").Value = 50
$ws.Range("N42").Value = 12000
$ws.Range("O42").Value = 12000
$ws.Range("P42").Value = 12000
$ws.Range("Q42").Value = '$/caja 18 kilos empedrada'
$ws.Range("R42").Value = 'Región de O''Higgins'
$ws.Range("S42").Value = 667
$ws.Range("T42").Value = 18

# Row 43
$ws.Range("D43").Value = 45079
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 30
$ws.Range("N43").Value = 10000
$ws.Range("O43").Value = 10000
$ws.Range("P43").Value = 10000
$ws.Range("Q43").Value = '$/caja 18 kilos empedrada'
$ws.Range("R43").Value = 'Región de O''Higgins'
$ws.Range("S43").Value = 556
$ws.Range("T43").Value = 18

# Row 44
$ws.Range("D44").Value = 45079
$ws.Range("L44").Value = 'Segunda'
$ws.Range("M44").Value = 20
$ws.Range("N44").Value = 9000
$ws.Range("O44").Value = 9000
$ws.Range("P44").Value = 9000
$ws.Range("Q44").Value = '$/caja 18 kilos empedrada'
$ws.Range("R44").Value = 'Región de O''Higgins'
$ws.Range("S44").Value = 500
$ws.Range("T44").Value = 18

# Row 45
$ws.Range("D45").Value = 45071
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 40
$ws.Range("N45").Value = 12000
$ws.Range("O45").Value = 12000
$ws.Range("P45").Value = 12000
$ws.Range("Q45").Value = '$/caja 18 kilos empedrada'
$ws.Range("R45").Value = 'Región de O''Higgins'
$ws.Range("S45").Value = 667
$ws.Range("T45").Value = 18

# Row 46
$ws.Range("D46").Value = 45071
$ws.Range("L46").Value = 'Segunda'
$ws.Range("M46").Value = 40
$ws.Range("N46").Value = 10000
$ws.Range("O46").Value = 10000
$ws.Range("P46").Value = 10000
$ws.Range("Q46").Value = '$/caja 18 kilos empedrada'
$ws.Range("R46").Value = 'Región de O''Higgins'
$ws.Range("S46").Value = 556
$ws.Range("T46").Value = 18

# Row 47
$ws.Range("D47").Value = 45070
$ws.Range("L47").Value = 'Primera'
$ws.Range("M47").Value = 60
$ws.Range("N47").Value = 10000
$ws.Range("O47").Value = 10000
$ws.Range("P47").Value = 10000
$ws.Range("Q47").Value = '$/caja 18 kilos empedrada'
$ws.Range("R47").Value = 'Región de O''Higgins'
$ws.Range("S47").Value = 556
$ws.Range("T47").Value = 18

# Row 48
$ws.Range("D48").Value = 45040
$ws.Range("L48").Value = 'Especial'
$ws.Range("M48").Value = 50
$ws.Range("N48").Value = 13000
$ws.Range("O48").Value = 13000
$ws.Range("P48").Value = 13000
$ws.Range("Q48").Value = '$/caja 18 kilos empedrada'
$ws.Range("R48").Value = 'Región de O''Higgins'
$ws.Range("S48").Value = 722
$ws.Range("T48").Value = 18

# Row 49
$ws.Range("D49").Value = 45040
$ws.Range("L49").Value = 'Primera'
$ws.Range("M49").Value = 40
$ws.Range("N49").Value = 12000
$ws.Range("O49").Value = 12000
$ws.Range("P49").Value = 12000
$ws.Range("Q49").Value = '$/caja 18 kilos empedrada'
$ws.Range("R49").Value = 'Región de O''Higgins'
$ws.Range("S49").Value = 667
$ws.Range("T49").Value = 18

# Row 50
$ws.Range("D50").Value = 45106
$ws.Range("L50").Value = 'Especial'
$ws.Range("M50").Value = 80
$ws.Range("N50").Value = 12000
$ws.Range("O50").Value = 12000
$ws.Range("P50").Value = 12000
$ws.Range("Q50").Value = '$/caja 18 kilos empedrada'
$ws.Range("R50").Value = 'Región del Maule'
$ws.Range("S50").Value = 667
$ws.Range("T50").Value = 18

# Row 51
$ws.Range("D51").Value = 45106
$ws.Range("L51").Value = 'Primera'
$ws.Range("M51").Value = 50
$ws.Range("N51").Value = 10000
$ws.Range("O51").Value = 10000
$ws.Range("P51").Value = 10000
$ws.Range("Q51").Value = '$/caja 18 kilos empedrada'
$ws.Range("R51").Value = 'Región del Maule'
$ws.Range("S51").Value = 556
$ws.Range("T51").Value = 18

# Row 52
$ws.Range("D52").Value = 45112
$ws.Range("L52").Value = 'Primera'
$ws.Range("M52").Value = 50
$ws.Range("N52").Value = 10000
$ws.Range("O52").Value = 10000
$ws.Range("P52").Value = 10000
$ws.Range("Q52").Value = '$/caja 18 kilos empedrada'
$ws.Range("R52").Value = 'Región de O''Higgins'
$ws.Range("S52").Value = 556
$ws.Range("T52").Value = 18

# Row 53
$ws.Range("D53").Value = 45112
$ws.Range("L53").Value = 'Segunda'
$ws.Range("M53").Value = 50
$ws.Range("N53").Value = 8000
$ws.Range("O53").Value = 8000
$ws.Range("P53").Value = 8000
$ws.Range("Q53").Value = '$/caja 18 kilos empedrada'
$ws.Range("R53").Value = 'Región de O''Higgins'
$ws.Range("S53").Value = 444
$ws.Range("T53").Value = 18

# Row 54
$ws.Range("D54").Value = 45089
$ws.Range("L54").Value = 'Especial'
$ws.Range("M54").Value = 60
$ws.Range("N54").Value = 11000
$ws.Range("O54").Value = 11000
$ws.Range("P54").Value = 11000
$ws.Range("Q54").Value = '$/caja 18 kilos empedrada'
$ws.Range("R54").Value = 'Región del Maule'
$ws.Range("S54").Value = 611
$ws.Range("T54").Value = 18

# Row 55
$ws.Range("D55").Value = 45089
$ws.Range("L55").Value = 'Primera'
$ws.Range("M55").Value = 50
$ws.Range("N55").Value = 9000
$ws.Range("O55").Value = 9000
$ws.Range("P55").Value = 9000
$ws.Range("Q55").Value = '$/caja 18 kilos empedrada'
$ws.Range("R55").Value = 'Región del Maule'
$ws.Range("S55").Value = 500
$ws.Range("T55").Value = 18

# Row 56
$ws.Range("D56").Value = 45089
$ws.Range("L56").Value = 'Segunda'
$ws.Range("M56").Value = 30
$ws.Range("N56").Value = 7000
$ws.Range("O56").Value = 7000
$ws.Range("P56").Value = 7000
$ws.Range("Q56").Value = '$/caja 18 kilos empedrada'
$ws.Range("R56").Value = 'Región del Maule'
$ws.Range("S56").Value = 389
$ws.Range("T56").Value = 18

# Row 57
$ws.Range("D57").Value = 45099
$ws.Range("L57").Value = 'Especial'
$ws.Range("M57").Value = 100
$ws.Range("N57").Value = 12000
$ws.Range("O57").Value = 12000
$ws.Range("P57").Value = 12000
$ws.Range("Q57").Value = '$/caja 18 kilos empedrada'
$ws.Range("R57").Value = 'Región del Maule'
$ws.Range("S57").Value = 667
$ws.Range("T57").Value = 18

# Row 58
$ws.Range("D58").Value = 45099
$ws.Range("L58").Value = 'Primera'
$ws.Range("M58").Value = 80
$ws.Range("N58").Value = 10000
$ws.Range("O58").Value = 10000
$ws.Range("P58").Value = 10000
$ws.Range("Q58").Value = '$/caja 18 kilos empedrada'
$ws.Range("R58").Value = 'Región del Maule'
$ws.Range("S58").Value = 556
$ws.Range("T58").Value = 18

# Row 59
$ws.Range("D59").Value = 45099
$ws.Range("L59").Value = 'Segunda'
$ws.Range("M59").Value = 80
$ws.Range("N59").Value = 8000
$ws.Range("O59").Value = 8000
$ws.Range("P59").Value = 8000
$ws.Range("Q59").Value = '$/caja 18 kilos empedrada'
$ws.Range("R59").Value = 'Región del Maule'
$ws.Range("S59").Value = 444
$ws.Range("T59").Value = 18
